$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (sampleCookie1)
$ws.Range("A2").Value = 45664.23578169493
$ws.Range("E2").Value = 45671.4441087963
$ws.Range("L2").Value = "Secure"
$ws.Range("N2").Value = 45664.23578169434
$ws.Range("O2").Value = "Other"

# Row 3 (sampleCookie2)
$ws.Range("A3").Value = 45664.23578172098
$ws.Range("E3").Value = 45694.4441087963
$ws.Range("L3").Value = "Secure"
$ws.Range("N3").Value = 45664.23578172042
$ws.Range("O3").Value = "Other"

# Row 4 (sampleCookie3)
$ws.Range("A4").Value = 45664.23578176276
$ws.Range("L4").Value = "Secure"
$ws.Range("N4").Value = 45664.23578176225
$ws.Range("O4").Value = "Other"
